$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures for cryptos.xlsx.
# For Price (column D) values that are plain decimal numbers, a leading
# apostrophe forces Excel to keep them as literal text (matching the
# original inline-string cell type) instead of silently converting them
# to numbers and losing formatting (e.g. trailing zeros, scientific notation).

$ws.Range("D2").Value = "27.796.04"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "1.866.13"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  +3.20%  "
$ws.Range("D5").Value = "'324.78"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").Value = "'0.4428"
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("D8").Value = "'0.3805"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "'0.07483"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").Value = "'0.8863"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "'21.74"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "1.891.14"
$ws.Range("E12").Value = "  -12.50%  "
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "'6.763"
$ws.Range("D15").Value = "'0.07239"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "'83.96"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "'0.000009158"
$ws.Range("E18").Value = "  +4.04%  "
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "27.793.78"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'5.330"
$ws.Range("E22").Value = "  +2.94%  "
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").Value = "'1.988"
$ws.Range("E24").Value = "  +5.32%  "
$ws.Range("D25").Value = "'158.80"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("D27").Value = "'5.337"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("D29").Value = "'117.89"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").Value = "'0.09085"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'3.123"
$ws.Range("E31").Value = "  +11.10%  "
$ws.Range("D32").Value = "'0.7802"
$ws.Range("E32").Value = "  +4.36%  "
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("D34").Value = "'4.582"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("D36").Value = "'1.154"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("D38").Value = "'0.05355"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").Value = "'2.865"
$ws.Range("E39").Value = "  +4.50%  "
$ws.Range("D40").Value = "'0.5211"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").Value = "'0.1697"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").Value = "'6.938"
$ws.Range("E42").Value = "  +6.96%  "
$ws.Range("D43").Value = "'8.686"
$ws.Range("E43").Value = "  +4.56%  "
$ws.Range("D44").Value = "'10.74"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").Value = "'109.78"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").Value = "'1.727"
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("D47").Value = "'0.4720"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "'0.06478"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").Value = "'1.912"
$ws.Range("E49").Value = "  +4.38%  "
$ws.Range("D50").Value = "'39.93"
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("E51").Value = "  +2.58%  "
